$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6899.7334
$ws.Range("J40").Value = 13748.75
$ws.Range("L40").Value = 13748.75
$ws.Range("N40").Value = -14098.75
$ws.Range("H92").Value = 444
$ws.Range("I92").Value = 444
$ws.Range("J92").Value = 444
$ws.Range("K92").Value = 444
$ws.Range("L92").Value = 444
$ws.Range("M92").Value = 804
$ws.Range("N92").Value = -2940
$ws.Range("H94").Value = 2075.2727
$ws.Range("I94").Value = 1572.75
$ws.Range("K94").Value = 1572.75
$ws.Range("M94").Value = -1121.75
$ws.Range("H96").Value = 564.1667
$ws.Range("I96").Value = 477.1
$ws.Range("K96").Value = 1431.3
$ws.Range("M96").Value = -58.30000000000018
$ws.Range("H100").Value = 7916.0435
$ws.Range("I100").Value = 1643.091
$ws.Range("K100").Value = 1643.091
$ws.Range("M100").Value = -1102.091
$ws.Range("H104").Value = 1004.1667
$ws.Range("I104").Value = 1004.1667
$ws.Range("K104").Value = 3012.5001
$ws.Range("M104").Value = -1265.5001
$ws.Range("H111").Value = 72236.07000000001
$ws.Range("I111").Value = 91636.82000000001
$ws.Range("K111").Value = 274910.46
$ws.Range("M111").Value = -271843.46
$ws.Range("H132").Value = 8278.882
$ws.Range("I132").Value = 1475.04
$ws.Range("K132").Value = 4425.12
$ws.Range("M132").Value = -1895.12
$ws.Range("H137").Value = 5182.7334
$ws.Range("I137").Value = 5269.8335
$ws.Range("J137").Value = 4834.3335
$ws.Range("K137").Value = 15809.5005
$ws.Range("L137").Value = 14503.0005
$ws.Range("M137").Value = -13259.5005
$ws.Range("N137").Value = -19603.0005
$ws.Range("H138").Value = 4537.8677
$ws.Range("J138").Value = 4874.1914
$ws.Range("L138").Value = 14622.5742
$ws.Range("N138").Value = -24902.5742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3977.3
$ws.Range("I32").Value = 3641.164
$ws.Range("K32").Value = 3641.164
$ws.Range("M32").Value = -3354.164
$ws.Range("H45").Value = 2105.8462
$ws.Range("I45").Value = 1852.5454
$ws.Range("J45").Value = 3499
$ws.Range("K45").Value = 1852.5454
$ws.Range("L45").Value = 3499
$ws.Range("M45").Value = -1475.5454
$ws.Range("N45").Value = -4253
$ws.Range("H74").Value = 1643.5238
$ws.Range("I74").Value = 1559.9412
$ws.Range("K74").Value = 1559.9412
$ws.Range("M74").Value = -685.9412
$ws.Range("H77").Value = 1643.5238
$ws.Range("I77").Value = 1559.9412
$ws.Range("K77").Value = 7799.706
$ws.Range("M77").Value = -3431.706
$ws.Range("H97").Value = 1266.5
$ws.Range("I97").Value = 1266.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1266.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -770.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 774826.6
$ws.Range("I86").Value = 1064524.8
$ws.Range("K86").Value = 1064524.8
$ws.Range("M86").Value = -1063401.8
$ws.Range("H89").Value = 774826.6
$ws.Range("I89").Value = 1064524.8
$ws.Range("K89").Value = 5322624
$ws.Range("M89").Value = -5317008
$ws.Range("H99").Value = 3460.2144
$ws.Range("I99").Value = 3695.8
$ws.Range("K99").Value = 3695.8
$ws.Range("M99").Value = -2197.8
$ws.Range("H107").Value = 669329.5600000001
$ws.Range("I107").Value = 2044.4
$ws.Range("K107").Value = 2044.4
$ws.Range("M107").Value = -124.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3202.889
$ws.Range("I58").Value = 6149.5
$ws.Range("J58").Value = 845.6
$ws.Range("K58").Value = 6149.5
$ws.Range("L58").Value = 845.6
$ws.Range("M58").Value = -5946.5
$ws.Range("N58").Value = -1251.6
$ws.Range("H136").Value = 3202.889
$ws.Range("I136").Value = 6149.5
$ws.Range("J136").Value = 845.6
$ws.Range("K136").Value = 18448.5
$ws.Range("L136").Value = 2536.8
$ws.Range("M136").Value = -15898.5
$ws.Range("N136").Value = -7636.8
$ws.Range("H141").Value = 328609.47
$ws.Range("I141").Value = 107796
$ws.Range("J141").Value = 368757.38
$ws.Range("K141").Value = 107796
$ws.Range("L141").Value = 368757.38
$ws.Range("M141").Value = -102616
$ws.Range("N141").Value = -379117.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2002.3846
$ws.Range("I63").Value = 1669.25
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 5007.75
$ws.Range("L63").Value = 18000
$ws.Range("M63").Value = -4258.75
$ws.Range("N63").Value = -19498
$ws.Range("H66").Value = 2002.3846
$ws.Range("I66").Value = 1669.25
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 15023.25
$ws.Range("L66").Value = 54000
$ws.Range("M66").Value = -11279.25
$ws.Range("N66").Value = -61488
$ws.Range("H107").Value = 89965.05
$ws.Range("I107").Value = 1131.8572
$ws.Range("J107").Value = 128829.56
$ws.Range("K107").Value = 3395.5716
$ws.Range("L107").Value = 386488.68
$ws.Range("M107").Value = -1475.5716
$ws.Range("N107").Value = -390328.68
$ws.Range("H140").Value = 3401.96
$ws.Range("I140").Value = 2134.158
$ws.Range("K140").Value = 6402.474
$ws.Range("M140").Value = -1222.474
$ws.Range("H141").Value = 14506.75
$ws.Range("I141").Value = 14506.75
$ws.Range("K141").Value = 43520.25
$ws.Range("M141").Value = -38340.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 65000
$ws.Range("J40").Value = 65000
$ws.Range("L40").Value = 65000
$ws.Range("N40").Value = -65302
$ws.Range("H97").Value = 1243
$ws.Range("I97").Value = 1342.125
$ws.Range("J97").Value = 978.6667
$ws.Range("K97").Value = 1342.125
$ws.Range("L97").Value = 978.6667
$ws.Range("M97").Value = -846.125
$ws.Range("N97").Value = -1970.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6536.125
$ws.Range("I7").Value = 6341.4287
$ws.Range("J7").Value = 7899
$ws.Range("K7").Value = 6341.4287
$ws.Range("L7").Value = 7899
$ws.Range("M7").Value = -6229.4287
$ws.Range("N7").Value = -8123
$ws.Range("H22").Value = 725
$ws.Range("J22").Value = 850.125
$ws.Range("L22").Value = 850.125
$ws.Range("N22").Value = -1440.125
$ws.Range("H27").Value = 725
$ws.Range("J27").Value = 850.125
$ws.Range("L27").Value = 850.125
$ws.Range("N27").Value = -1064.125
$ws.Range("H68").Value = 168466.67
$ws.Range("I68").Value = 966.6667
$ws.Range("K68").Value = 966.6667
$ws.Range("M68").Value = -217.6667
$ws.Range("H71").Value = 168466.67
$ws.Range("I71").Value = 966.6667
$ws.Range("K71").Value = 4833.3335
$ws.Range("M71").Value = -1089.3335
$ws.Range("H82").Value = 1197
$ws.Range("I82").Value = 1339.625
$ws.Range("K82").Value = 1339.625
$ws.Range("M82").Value = -978.625
$ws.Range("H85").Value = 1197
$ws.Range("I85").Value = 1339.625
$ws.Range("K85").Value = 1339.625
$ws.Range("M85").Value = -91.625
$ws.Range("H122").Value = 6258
$ws.Range("J122").Value = 6779
$ws.Range("L122").Value = 20337
$ws.Range("N122").Value = -25237
$ws.Range("H126").Value = 6536.125
$ws.Range("I126").Value = 6341.4287
$ws.Range("J126").Value = 7899
$ws.Range("K126").Value = 19024.2861
$ws.Range("L126").Value = 23697
$ws.Range("M126").Value = -16554.2861
$ws.Range("N126").Value = -28637
$ws.Range("H132").Value = 5106.744
$ws.Range("I132").Value = 4654.5864
$ws.Range("J132").Value = 6043.357
$ws.Range("K132").Value = 13963.7592
$ws.Range("L132").Value = 18130.071
$ws.Range("M132").Value = -11433.7592
$ws.Range("N132").Value = -23190.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 981.94116
$ws.Range("J107").Value = 660
$ws.Range("L107").Value = 1980
$ws.Range("N107").Value = -5820
$ws.Range("H113").Value = 187.625
$ws.Range("I113").Value = 187.28572
$ws.Range("K113").Value = 561.85716
$ws.Range("M113").Value = 1608.14284
$ws.Range("H132").Value = 25784.91
$ws.Range("I132").Value = 1310.2916
$ws.Range("J132").Value = 53755.906
$ws.Range("K132").Value = 3930.8748
$ws.Range("L132").Value = 161267.718
$ws.Range("M132").Value = -1400.8748
$ws.Range("N132").Value = -166327.718

# Clear cells that should no longer contain a value
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N97").ClearContents()
